$wb = $excel.ActiveWorkbook

# --- Repayment Schedule sheet: insert a new blank column before column N (14) ---
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = 9.17

# --- Switch the active sheet/tab to "Repayment Schedule" and move the selection ---
$ws.Activate()
[void]$ws.Range("R8").Select()
